# Update the "Förändrad" (Changed) date column (C) for all data rows.
# All cells in C2:C338 currently hold the Excel serial date 45177
# (2023-09-08) and must be bumped to 45178 (2023-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$ws.Range("C2:C$lastRow").Value = 45178
